# Apply updated market-price derived values (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR Leve-profit tables, per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3710.2307
$ws.Range("I70").Value = 1674.75
$ws.Range("K70").Value = 5024.25
$ws.Range("M70").Value = -4754.25
$ws.Range("H73").Value = 3710.2307
$ws.Range("I73").Value = 1674.75
$ws.Range("K73").Value = 5024.25
$ws.Range("M73").Value = -4088.25
$ws.Range("H100").Value = 1514.7778
$ws.Range("I100").Value = 947.8570999999999
$ws.Range("J100").Value = 3499
$ws.Range("K100").Value = 947.8570999999999
$ws.Range("L100").Value = 3499
$ws.Range("M100").Value = -406.8570999999999
$ws.Range("N100").Value = -4581
$ws.Range("H116").Value = 4400.3335
$ws.Range("I116").Value = 4101
$ws.Range("J116").Value = 4550
$ws.Range("K116").Value = 4101
$ws.Range("L116").Value = 4550
$ws.Range("M116").Value = -659
$ws.Range("N116").Value = -11434
$ws.Range("H132").Value = 3221.7693
$ws.Range("I132").Value = 3221.7693
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9665.3079
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7135.3079
$ws.Range("N132").ClearContents()
$ws.Range("H138").Value = 3378.0637
$ws.Range("J138").Value = 2665.7368
$ws.Range("L138").Value = 7997.2104
$ws.Range("N138").Value = -18277.2104

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1620
$ws.Range("I61").Value = 1517.238
$ws.Range("K61").Value = 1517.238
$ws.Range("M61").Value = -1305.238
$ws.Range("H102").Value = 1630.5
$ws.Range("I102").Value = 1519.8
$ws.Range("K102").Value = 1519.8
$ws.Range("M102").Value = 102.2
$ws.Range("H122").Value = 1670791.5
$ws.Range("I122").Value = 5002374.5
$ws.Range("K122").Value = 15007123.5
$ws.Range("M122").Value = -15004673.5
$ws.Range("H132").Value = 2552.6155
$ws.Range("I132").Value = 2538.45
$ws.Range("J132").Value = 2599.8333
$ws.Range("K132").Value = 7615.349999999999
$ws.Range("L132").Value = 7799.499899999999
$ws.Range("M132").Value = -5085.349999999999
$ws.Range("N132").Value = -12859.4999
$ws.Range("H136").Value = 1620
$ws.Range("I136").Value = 1517.238
$ws.Range("K136").Value = 4551.714
$ws.Range("M136").Value = -2001.714

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 9500
$ws.Range("I24").Value = 9500
$ws.Range("K24").Value = 9500
$ws.Range("M24").Value = -9265
$ws.Range("H86").Value = 2989.4
$ws.Range("I86").Value = 3043.7778
$ws.Range("J86").Value = 2500
$ws.Range("K86").Value = 3043.7778
$ws.Range("L86").Value = 2500
$ws.Range("M86").Value = -1920.7778
$ws.Range("N86").Value = -4746
$ws.Range("H89").Value = 2989.4
$ws.Range("I89").Value = 3043.7778
$ws.Range("J89").Value = 2500
$ws.Range("K89").Value = 15218.889
$ws.Range("L89").Value = 12500
$ws.Range("M89").Value = -9602.888999999999
$ws.Range("N89").Value = -23732
$ws.Range("H134").Value = 2233.7917
$ws.Range("I134").Value = 2034.6111
$ws.Range("K134").Value = 6103.8333
$ws.Range("M134").Value = -3568.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4356.5557
$ws.Range("I31").Value = 2797.75
$ws.Range("K31").Value = 2797.75
$ws.Range("M31").Value = -2502.75
$ws.Range("H34").Value = 4356.5557
$ws.Range("I34").Value = 2797.75
$ws.Range("K34").Value = 2797.75
$ws.Range("M34").Value = -2595.75
$ws.Range("H58").Value = 3237.04
$ws.Range("I58").Value = 1163.8182
$ws.Range("K58").Value = 1163.8182
$ws.Range("M58").Value = -960.8181999999999
$ws.Range("H107").Value = 905
$ws.Range("I107").Value = 905
$ws.Range("K107").Value = 905
$ws.Range("M107").Value = 1015
$ws.Range("H132").Value = 1932.625
$ws.Range("I132").Value = 2063.6296
$ws.Range("K132").Value = 6190.888800000001
$ws.Range("M132").Value = -3660.888800000001
$ws.Range("H134").Value = 2237.9167
$ws.Range("I134").Value = 1026.8572
$ws.Range("J134").Value = 3933.4
$ws.Range("K134").Value = 3080.5716
$ws.Range("L134").Value = 11800.2
$ws.Range("M134").Value = -545.5715999999998
$ws.Range("N134").Value = -16870.2
$ws.Range("H136").Value = 3237.04
$ws.Range("I136").Value = 1163.8182
$ws.Range("K136").Value = 3491.4546
$ws.Range("M136").Value = -941.4546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 290.5
$ws.Range("I9").Value = 81
$ws.Range("J9").Value = 500
$ws.Range("K9").Value = 243
$ws.Range("L9").Value = 1500
$ws.Range("M9").Value = -19
$ws.Range("N9").Value = -1948
$ws.Range("H131").Value = 1260.9546
$ws.Range("J131").Value = 1478.2778
$ws.Range("L131").Value = 4434.8334
$ws.Range("N131").Value = -14514.8334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 1620
$ws.Range("I29").Value = 1775
$ws.Range("K29").Value = 1775
$ws.Range("M29").Value = -1485
$ws.Range("H113").Value = 1100
$ws.Range("I113").Value = 1100
$ws.Range("K113").Value = 1100
$ws.Range("M113").Value = 1070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3016.9167
$ws.Range("I68").Value = 3022.75
$ws.Range("J68").Value = 3014
$ws.Range("K68").Value = 3022.75
$ws.Range("L68").Value = 3014
$ws.Range("M68").Value = -2273.75
$ws.Range("N68").Value = -4512
$ws.Range("H71").Value = 3016.9167
$ws.Range("I71").Value = 3022.75
$ws.Range("J71").Value = 3014
$ws.Range("K71").Value = 15113.75
$ws.Range("L71").Value = 15070
$ws.Range("M71").Value = -11369.75
$ws.Range("N71").Value = -22558
$ws.Range("H100").Value = 1436.625
$ws.Range("I100").Value = 1082.1666
$ws.Range("K100").Value = 1082.1666
$ws.Range("M100").Value = -541.1666
$ws.Range("H136").Value = 1780.341
$ws.Range("I136").Value = 1714.6578
$ws.Range("K136").Value = 5143.9734
$ws.Range("M136").Value = -2593.9734

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 10929
$ws.Range("I63").Value = 10929
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 10929
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -10305
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 10929
$ws.Range("I66").Value = 10929
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 32787
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -29667
$ws.Range("N66").ClearContents()
$ws.Range("H81").Value = 5761.75
$ws.Range("J81").Value = 750
$ws.Range("L81").Value = 1500
$ws.Range("N81").Value = -3622
$ws.Range("H84").Value = 5761.75
$ws.Range("J84").Value = 750
$ws.Range("L84").Value = 7500
$ws.Range("N84").Value = -18108
$ws.Range("H96").Value = 1290.7
$ws.Range("I96").Value = 1288.8572
$ws.Range("K96").Value = 1288.8572
$ws.Range("M96").Value = 84.14280000000008
$ws.Range("H100").Value = 2224.5557
$ws.Range("I100").Value = 2003.1428
$ws.Range("K100").Value = 4006.2856
$ws.Range("M100").Value = -3465.2856
$ws.Range("H122").Value = 2049.8572
$ws.Range("I122").Value = 2049.8572
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6149.571599999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3699.571599999999
$ws.Range("N122").ClearContents()
